$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3668.1428
$ws.Range("I86").Value = 3094.5386
$ws.Range("K86").Value = 3094.5386
$ws.Range("M86").Value = -1971.5386
$ws.Range("H89").Value = 3668.1428
$ws.Range("I89").Value = 3094.5386
$ws.Range("K89").Value = 15472.693
$ws.Range("M89").Value = -9856.692999999999
$ws.Range("H120").Value = 49687.5
$ws.Range("J120").Value = 49687.5
$ws.Range("L120").Value = 49687.5
$ws.Range("N120").Value = -59363.5
$ws.Range("H128").Value = 37998.8
$ws.Range("J128").Value = 37998.8
$ws.Range("L128").Value = 37998.8
$ws.Range("N128").Value = -47958.8
$ws.Range("H137").Value = 4372.756
$ws.Range("J137").Value = 5138.879
$ws.Range("L137").Value = 15416.637
$ws.Range("N137").Value = -20516.637

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3527
$ws.Range("I61").Value = 2389
$ws.Range("J61").Value = 4125.9473
$ws.Range("K61").Value = 2389
$ws.Range("L61").Value = 4125.9473
$ws.Range("M61").Value = -2177
$ws.Range("N61").Value = -4549.9473
$ws.Range("H107").Value = 45224
$ws.Range("J107").Value = 45224
$ws.Range("L107").Value = 45224
$ws.Range("N107").Value = -52904
$ws.Range("H109").Value = 39924.4
$ws.Range("J109").Value = 39924.4
$ws.Range("L109").Value = 39924.4
$ws.Range("N109").Value = -42698.4
$ws.Range("H111").Value = 45328
$ws.Range("J111").Value = 45328
$ws.Range("L111").Value = 45328
$ws.Range("N111").Value = -53508
$ws.Range("H117").Value = 42847.8
$ws.Range("J117").Value = 42847.8
$ws.Range("L117").Value = 42847.8
$ws.Range("N117").Value = -52025.8
$ws.Range("H124").Value = 20476.334
$ws.Range("J124").Value = 20476.334
$ws.Range("L124").Value = 20476.334
$ws.Range("N124").Value = -30296.334
$ws.Range("H125").Value = 50707
$ws.Range("J125").Value = 50707
$ws.Range("L125").Value = 50707
$ws.Range("N125").Value = -60547
$ws.Range("H128").Value = 48423.668
$ws.Range("J128").Value = 48423.668
$ws.Range("L128").Value = 48423.668
$ws.Range("N128").Value = -58383.668
$ws.Range("H136").Value = 3527
$ws.Range("I136").Value = 2389
$ws.Range("J136").Value = 4125.9473
$ws.Range("K136").Value = 7167
$ws.Range("L136").Value = 12377.8419
$ws.Range("M136").Value = -4617
$ws.Range("N136").Value = -17477.8419

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H108").Value = 47684
$ws.Range("J108").Value = 47684
$ws.Range("L108").Value = 47684
$ws.Range("N108").Value = -55364
$ws.Range("H110").Value = 48080.5
$ws.Range("J110").Value = 48080.5
$ws.Range("L110").Value = 48080.5
$ws.Range("N110").Value = -56260.5
$ws.Range("H112").Value = 46361.25
$ws.Range("J112").Value = 46361.25
$ws.Range("L112").Value = 46361.25
$ws.Range("N112").Value = -49315.25
$ws.Range("H124").Value = 47992
$ws.Range("J124").Value = 47992
$ws.Range("L124").Value = 47992
$ws.Range("N124").Value = -57812

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 45459.6
$ws.Range("J20").Value = 45459.6
$ws.Range("L20").Value = 45459.6
$ws.Range("N20").Value = -45931.6
$ws.Range("H30").Value = 45459.6
$ws.Range("J30").Value = 45459.6
$ws.Range("L30").Value = 45459.6
$ws.Range("N30").Value = -45641.6
$ws.Range("H31").Value = 159695.4
$ws.Range("I31").Value = 1797.7084
$ws.Range("J31").Value = 217112.75
$ws.Range("K31").Value = 1797.7084
$ws.Range("L31").Value = 217112.75
$ws.Range("M31").Value = -1502.7084
$ws.Range("N31").Value = -217702.75
$ws.Range("H34").Value = 159695.4
$ws.Range("I34").Value = 1797.7084
$ws.Range("J34").Value = 217112.75
$ws.Range("K34").Value = 1797.7084
$ws.Range("L34").Value = 217112.75
$ws.Range("M34").Value = -1595.7084
$ws.Range("N34").Value = -217516.75
$ws.Range("H116").Value = 49368.5
$ws.Range("J116").Value = 49368.5
$ws.Range("L116").Value = 49368.5
$ws.Range("N116").Value = -58546.5
$ws.Range("H125").Value = 16856.8
$ws.Range("J125").Value = 16856.8
$ws.Range("L125").Value = 16856.8
$ws.Range("N125").Value = -21776.8
$ws.Range("H128").Value = 45459.6
$ws.Range("J128").Value = 45459.6
$ws.Range("L128").Value = 45459.6
$ws.Range("N128").Value = -55419.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1388
$ws.Range("I17").Value = 41
$ws.Range("J17").Value = 1724.75
$ws.Range("K17").Value = 123
$ws.Range("L17").Value = 5174.25
$ws.Range("M17").Value = 46
$ws.Range("N17").Value = -5512.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4795.5557
$ws.Range("I70").Value = 4831.515
$ws.Range("J70").Value = 4400
$ws.Range("K70").Value = 4831.515
$ws.Range("L70").Value = 4400
$ws.Range("M70").Value = -4561.515
$ws.Range("N70").Value = -4940
$ws.Range("H73").Value = 4795.5557
$ws.Range("I73").Value = 4831.515
$ws.Range("J73").Value = 4400
$ws.Range("K73").Value = 4831.515
$ws.Range("L73").Value = 4400
$ws.Range("M73").Value = -3895.515
$ws.Range("N73").Value = -6272
$ws.Range("H80").Value = 337594
$ws.Range("I80").Value = 458591.8
$ws.Range("K80").Value = 458591.8
$ws.Range("M80").Value = -457593.8
$ws.Range("H83").Value = 337594
$ws.Range("I83").Value = 458591.8
$ws.Range("K83").Value = 2292959
$ws.Range("M83").Value = -2287967
$ws.Range("H104").Value = 33797.2
$ws.Range("J104").Value = 33797.2
$ws.Range("L104").Value = 33797.2
$ws.Range("N104").Value = -40785.2
$ws.Range("H105").Value = 38519.855
$ws.Range("J105").Value = 38519.855
$ws.Range("L105").Value = 38519.855
$ws.Range("N105").Value = -45507.855
$ws.Range("H110").Value = 34628.5
$ws.Range("J110").Value = 34628.5
$ws.Range("L110").Value = 34628.5
$ws.Range("N110").Value = -42808.5
$ws.Range("H118").Value = 38298
$ws.Range("J118").Value = 38298
$ws.Range("L118").Value = 38298
$ws.Range("N118").Value = -41612
$ws.Range("H120").Value = 39309
$ws.Range("J120").Value = 39309
$ws.Range("L120").Value = 39309
$ws.Range("N120").Value = -48985
$ws.Range("H127").Value = 37330.668
$ws.Range("J127").Value = 37330.668
$ws.Range("L127").Value = 37330.668
$ws.Range("N127").Value = -47250.668
$ws.Range("H131").Value = 38986
$ws.Range("J131").Value = 38986
$ws.Range("L131").Value = 38986
$ws.Range("N131").Value = -49066

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 31993.666
$ws.Range("J106").Value = 31993.666
$ws.Range("L106").Value = 31993.666
$ws.Range("N106").Value = -34517.666
$ws.Range("H109").Value = 35277
$ws.Range("J109").Value = 35277
$ws.Range("L109").Value = 35277
$ws.Range("N109").Value = -38051
$ws.Range("H110").Value = 33756
$ws.Range("J110").Value = 33756
$ws.Range("L110").Value = 33756
$ws.Range("N110").Value = -41936
$ws.Range("H111").Value = 43940.5
$ws.Range("J111").Value = 43940.5
$ws.Range("L111").Value = 43940.5
$ws.Range("N111").Value = -52120.5
$ws.Range("H112").Value = 25758.572
$ws.Range("J112").Value = 27551.666
$ws.Range("L112").Value = 27551.666
$ws.Range("N112").Value = -30505.666
$ws.Range("H117").Value = 45380
$ws.Range("J117").Value = 45380
$ws.Range("L117").Value = 45380
$ws.Range("N117").Value = -54558
$ws.Range("H119").Value = 47412
$ws.Range("J119").Value = 47412
$ws.Range("L119").Value = 47412
$ws.Range("N119").Value = -57088
$ws.Range("H121").Value = 21259
$ws.Range("J121").Value = 21259
$ws.Range("L121").Value = 21259
$ws.Range("N121").Value = -24753
$ws.Range("H125").Value = 49707
$ws.Range("J125").Value = 49707
$ws.Range("L125").Value = 49707
$ws.Range("N125").Value = -59547

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 31816
$ws.Range("J109").Value = 31816
$ws.Range("L109").Value = 31816
$ws.Range("N109").Value = -34590
$ws.Range("H115").Value = 31713.8
$ws.Range("J115").Value = 31713.8
$ws.Range("L115").Value = 31713.8
$ws.Range("N115").Value = -34847.8
$ws.Range("H116").Value = 26556.5
$ws.Range("J116").Value = 26556.5
$ws.Range("L116").Value = 26556.5
$ws.Range("N116").Value = -35734.5
$ws.Range("H117").Value = 44704.5
$ws.Range("J117").Value = 44704.5
$ws.Range("L117").Value = 44704.5
$ws.Range("N117").Value = -53882.5
$ws.Range("H118").Value = 24670.4
$ws.Range("J118").Value = 27088
$ws.Range("L118").Value = 27088
$ws.Range("N118").Value = -30402
$ws.Range("H121").Value = 44412
$ws.Range("J121").Value = 44412
$ws.Range("L121").Value = 44412
$ws.Range("N121").Value = -47906
